$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for the new summary columns (K1:O1)
$ws.Range("K1").Value = "Sum"
$ws.Range("L1").Value = "Average"
$ws.Range("M1").Value = "Count"
$ws.Range("N1").Value = "CountA"
$ws.Range("O1").Value = "CountIF"

# Row 3 formulas: SUM / AVERAGE / COUNT / COUNTA demo
$ws.Range("K3").Formula = "=SUM(C3:H3)"
$ws.Range("L3").Formula = "=AVERAGE(C3:H3)"
$ws.Range("M3").Formula = "=COUNT(A3:H3)"
$ws.Range("N3").Formula = "=COUNTA(A3:H3)"

# Additional SUM variants in column K for rows 4-6
$ws.Range("K4").Formula = "=SUM(C4:E4,E6:H6)"
$ws.Range("K6").Formula = "=SUM(C6,D6,E6,F6,G6,H6)"

# Highlight K4 with a yellow fill
$ws.Range("K4").Interior.Color = 65535

# K5 uses the "+" addition form, then gets a 2-decimal number format
$ws.Range("K5").Formula = "=SUM(C5+D5+E5+F5+G5+H5)"
$ws.Range("K5").NumberFormat = "#,##0.00"

# Column widths: J (spacer) and K (values)
$ws.Columns.Item(10).ColumnWidth = 2.5
$ws.Columns.Item(11).ColumnWidth = 10.75

# Move selection to N3 to match the saved selection state
$ws.Range("N3").Select()
